$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Raw Data"
$ws.Name = "Raw Data"

# Copy the existing header format (bold, centered, bordered) onto the new header cells
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)

# Add new header columns H1 and I1
$ws.Range("H1").Value = "timestamp"
$ws.Range("I1").Value = "hour_slot"
